# Nasar - 27th Dec
# Updates the OneCSAppiOS object-repository sheet:
#   1. ORDER_LIST_TAB's value is renamed from "Order list" to "order-list-tab-button"
#   2. A batch of new object-repository entries is appended (Activity tab, pay-money-in
#      button, order status labels, order detail screen fields, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 155 (ORDER_LIST_TAB) value column changes.
$ws.Cells.Item(155, 2).Value = 'order-list-tab-button'

# 2) Append the new rows (columns: A=ELEMENT_KEY, B=ELEMENT_VALUE, E=VALUE_TYPE).
$newRows = @(
    @('ACTIVITY_TAB', 'activity-tab-button', 'device-accessibilityid'),
    @('PAY_MONEY_IN_BTN_ACC_MANAGEMENT', 'Pay money in', 'device-accessibilityid'),
    @('ORDER_TYPE_CANCELLED', '//XCUIElementTypeStaticText[@label="CANCELLED"]', 'device-xpath'),
    @('ORDER_TYPE_COMPLETE', '//XCUIElementTypeStaticText[@label="COMPLETE"]', 'device-xpath'),
    @('ORDER_TYPE_PROCESSING', '//XCUIElementTypeStaticText[@label="PROCESSING"]', 'device-xpath'),
    @('ORDER_TYPE_DEALT', '//XCUIElementTypeStaticText[@label="DEALT"]', 'device-xpath'),
    @('ORDER_TYPE_REJECTED', '//XCUIElementTypeStaticText[@label="REJECTED"]', 'device-xpath'),
    @('ORDER_TYPE_EXPIRED', '//XCUIElementTypeStaticText[@label="EXPIRED"]', 'device-xpath'),
    @('CLOSE_ORDER_DETAILS_BTN', 'close-order-details-button', 'device-accessibilityid'),
    @('ORDER_DETAILS_TITLE', 'Order detail', 'device-accessibilityid'),
    @('ORDER_DETAIL_STATUS_FIELD', 'Status', 'device-accessibilityid'),
    @('ORDER_DETAIL_TRADE_TYPE_FIELD', 'Trade type', 'device-accessibilityid'),
    @('ORDER_DETAIL_DATE_FIELD', '(//XCUIElementTypeStaticText[@name="Date"])[9]', 'device-xpath'),
    @('ORDER_DETAIL_SETTLEMENT_DATE_FIELD', 'Settlement date', 'device-accessibilityid'),
    @('ORDER_DETAIL_PRICE_FIELD', 'Price', 'device-accessibilityid'),
    @('ORDER_DETAIL_QUANTITY_FIELD', 'Quantity', 'device-accessibilityid'),
    @('ORDER_DETAIL_CONSIDERATION_FIELD', 'Consideration', 'device-accessibilityid'),
    @('ORDER_DETAIL_OTHER_FIELD', 'Other', 'device-accessibilityid'),
    @('ORDER_DETAIL_DEALING_FEE_FIELD', 'Dealing fee', 'device-accessibilityid'),
    @('ORDER_DETAIL_ESTIMATED_VALUE_FIELD', 'Estimated value', 'device-accessibilityid'),
    @('ORDER_DETAIL_CANCEL_ORDER_BTN', 'cancel-order-button', 'device-accessibilityid')
)

$startRow = 206
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $entry = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 5).Value = $entry[2]
}

# Move the visible selection to match the author's final cursor position.
$ws.Range("B229").Select()
